# Automatically update sets from .txt files
# Appends three new LEGO set rows (15-17) to Sheet1, mirroring the
# existing table layout (ID_Set, Nom_Set, nbPieces, Collection, Image_URL,
# URL_Amazon, URL_Lego, URL_Auchan, URL_Leclerc, URL_Carrefour, URL_Fnac).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (ID_Set) and C (nbPieces) hold numeric-looking values that are
# stored as text in this workbook, so force text formatting before writing
# them (keeps "43230" etc. from being reinterpreted as numbers).
$ws.Range("A15:A17").NumberFormat = "@"
$ws.Range("C15:C17").NumberFormat = "@"

# Row 15 - 43230 La caméra Hommage à Walt Disney
$ws.Range("A15").Value = "43230"
$ws.Range("B15").Value = "La caméra Hommage à Walt Disney"
$ws.Range("C15").Value = "811"
$ws.Range("D15").Value = "Disney™"
$ws.Range("E15").Value = "https://www.lego.com/cdn/cs/set/assets/blta0d2ef903df1c30c/43230.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Range("F15").Value = "https://amzn.eu/d/b1S8pul"
$ws.Range("G15").Value = "https://www.lego.com/fr-fr/product/43230"
$ws.Range("H15").Value = "https://www.auchan.fr/lego-lego-disney-43230-la-camera-hommage-a-walt-disney-maquette-pour-adultes-avec-mickey-et-minnie-mouse/pr-C1718290"

# Row 16 - 10368 Le chrysanthème
$ws.Range("A16").Value = "10368"
$ws.Range("B16").Value = "Le chrysanthème"
$ws.Range("C16").Value = "278"
$ws.Range("D16").Value = "The Botanical Collection"
$ws.Range("E16").Value = "https://www.lego.com/cdn/cs/set/assets/bltdbc3129b50f61480/10368_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Range("F16").Value = "https://amzn.eu/d/hoKJjfA"
$ws.Range("G16").Value = "https://www.lego.com/fr-fr/product/10368"
$ws.Range("H16").Value = "https://www.auchan.fr/lego-icons-10368-chrysantheme-collection-botanique/pr-C1802539"
$ws.Range("I16").Value = "https://www.e.leclerc/fp/lego-icons-10368-le-chrysantheme-set-de-construction-5702017719689"
$ws.Range("J16").Value = "https://www.carrefour.fr/p/lego-icons-le-chrysantheme-10368-lego-5702017719689"

# Row 17 - 43257 Angel
$ws.Range("A17").Value = "43257"
$ws.Range("B17").Value = "Angel"
$ws.Range("C17").Value = "784"
$ws.Range("D17").Value = "Disney™"
$ws.Range("E17").Value = "https://www.lego.com/cdn/cs/set/assets/blt56c61562d64dc2e4/43257_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Range("F17").Value = "https://amzn.eu/d/eRPMY6r"
$ws.Range("G17").Value = "https://www.lego.com/fr-fr/product/43257"
$ws.Range("H17").Value = "https://www.auchan.fr/lego-disney-43257-angel-stitch/pr-C1836201"
$ws.Range("I17").Value = "https://www.e.leclerc/fp/lego-disney-angel-jouet-de-construction-lilo-et-stitch-pour-filles-et-garcons-43257-5702017813967"
$ws.Range("J17").Value = "https://www.carrefour.fr/p/lego-angel-43257-lego-5702017813967"
